$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1598.6765
$ws.Range("J17").Value = 1311.8788
$ws.Range("L17").Value = 3935.6364
$ws.Range("N17").Value = -4271.636399999999
$ws.Range("H41").Value = 612.41174
$ws.Range("I41").Value = 286.66666
$ws.Range("J41").Value = 682.2143
$ws.Range("K41").Value = 286.66666
$ws.Range("L41").Value = 682.2143
$ws.Range("M41").Value = 153.33334
$ws.Range("N41").Value = -1562.2143
$ws.Range("H74").Value = 4333.3335
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
$ws.Range("H77").Value = 4333.3335
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
$ws.Range("H92").Value = 2052591.4
$ws.Range("I92").Value = 2462910.2
$ws.Range("K92").Value = 2462910.2
$ws.Range("M92").Value = -2461662.2
$ws.Range("H94").Value = 3322.1667
$ws.Range("I94").Value = 2986.6
$ws.Range("K94").Value = 2986.6
$ws.Range("M94").Value = -2535.6
$ws.Range("H103").Value = 1167.44
$ws.Range("I103").Value = 999.4545000000001
$ws.Range("K103").Value = 2998.3635
$ws.Range("M103").Value = -2412.3635
$ws.Range("H125").Value = 1118.2142
$ws.Range("J125").Value = 1037.6666
$ws.Range("L125").Value = 9338.999400000001
$ws.Range("N125").Value = -14258.9994
$ws.Range("H129").Value = 871.80304
$ws.Range("I129").Value = 679.4
$ws.Range("K129").Value = 2038.2
$ws.Range("M129").Value = 2961.8
$ws.Range("H137").Value = 1995.1177
$ws.Range("I137").Value = 1666.4286
$ws.Range("J137").Value = 2225.2
$ws.Range("K137").Value = 4999.2858
$ws.Range("L137").Value = 6675.599999999999
$ws.Range("M137").Value = -2449.2858
$ws.Range("N137").Value = -11775.6
$ws.Range("H138").Value = 1733.3485
$ws.Range("I138").Value = 1216.8334
$ws.Range("J138").Value = 2163.7778
$ws.Range("K138").Value = 3650.5002
$ws.Range("L138").Value = 6491.3334
$ws.Range("M138").Value = 1489.4998
$ws.Range("N138").Value = -16771.3334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3052.443
$ws.Range("I32").Value = 1864.1177
$ws.Range("K32").Value = 1864.1177
$ws.Range("M32").Value = -1577.1177
$ws.Range("H74").Value = 2326.5881
$ws.Range("I74").Value = 1759.4
$ws.Range("K74").Value = 1759.4
$ws.Range("M74").Value = -885.4000000000001
$ws.Range("H77").Value = 2326.5881
$ws.Range("I77").Value = 1759.4
$ws.Range("K77").Value = 8797
$ws.Range("M77").Value = -4429
$ws.Range("H97").Value = 1057.6154
$ws.Range("I97").Value = 804.1111
$ws.Range("J97").Value = 1628
$ws.Range("K97").Value = 804.1111
$ws.Range("L97").Value = 1628
$ws.Range("M97").Value = -308.1111
$ws.Range("N97").Value = -2620
$ws.Range("H110").Value = 1544.9565
$ws.Range("I110").Value = 1160.0454
$ws.Range("K110").Value = 1160.0454
$ws.Range("M110").Value = 884.9546
$ws.Range("H125").Value = 49991.625
$ws.Range("J125").Value = 49991.625
$ws.Range("L125").Value = 49991.625
$ws.Range("N125").Value = -59831.625
$ws.Range("H132").Value = 1877.4839
$ws.Range("I132").Value = 1105.1428
$ws.Range("K132").Value = 3315.4284
$ws.Range("M132").Value = -785.4284000000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 199.75
$ws.Range("I22").Value = 199.75
$ws.Range("K22").Value = 199.75
$ws.Range("M22").Value = -26.75
$ws.Range("H86").Value = 89453.52
$ws.Range("I86").Value = 3597.6
$ws.Range("J86").Value = 155496.53
$ws.Range("K86").Value = 3597.6
$ws.Range("L86").Value = 155496.53
$ws.Range("M86").Value = -2474.6
$ws.Range("N86").Value = -157742.53
$ws.Range("H89").Value = 89453.52
$ws.Range("I89").Value = 3597.6
$ws.Range("J89").Value = 155496.53
$ws.Range("K89").Value = 17988
$ws.Range("L89").Value = 777482.65
$ws.Range("M89").Value = -12372
$ws.Range("N89").Value = -788714.65
$ws.Range("H94").Value = 703.9375
$ws.Range("I94").Value = 772.6667
$ws.Range("J94").Value = 497.75
$ws.Range("K94").Value = 772.6667
$ws.Range("L94").Value = 497.75
$ws.Range("M94").Value = -321.6667
$ws.Range("N94").Value = -1399.75
$ws.Range("H134").Value = 2522.4634
$ws.Range("I134").Value = 2387
$ws.Range("K134").Value = 7161
$ws.Range("M134").Value = -4626

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H31").Value = 2483.7273
$ws.Range("I31").Value = 2270.2222
$ws.Range("J31").Value = 2631.5386
$ws.Range("K31").Value = 2270.2222
$ws.Range("L31").Value = 2631.5386
$ws.Range("M31").Value = -1975.2222
$ws.Range("N31").Value = -3221.5386
$ws.Range("H34").Value = 2483.7273
$ws.Range("I34").Value = 2270.2222
$ws.Range("J34").Value = 2631.5386
$ws.Range("K34").Value = 2270.2222
$ws.Range("L34").Value = 2631.5386
$ws.Range("M34").Value = -2068.2222
$ws.Range("N34").Value = -3035.5386
$ws.Range("H47").Value = 10280
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 10280
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 10280
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -11412
$ws.Range("H58").Value = 1740758.1
$ws.Range("I58").Value = 2558592
$ws.Range("J58").Value = 2860.875
$ws.Range("K58").Value = 2558592
$ws.Range("L58").Value = 2860.875
$ws.Range("M58").Value = -2558389
$ws.Range("N58").Value = -3266.875
$ws.Range("H94").Value = 1167.875
$ws.Range("I94").Value = 1275
$ws.Range("K94").Value = 1275
$ws.Range("M94").Value = -824
$ws.Range("H134").Value = 1100.5814
$ws.Range("I134").Value = 1078.75
$ws.Range("J134").Value = 1212.8572
$ws.Range("K134").Value = 3236.25
$ws.Range("L134").Value = 3638.5716
$ws.Range("M134").Value = -701.25
$ws.Range("N134").Value = -8708.571599999999
$ws.Range("H136").Value = 1740758.1
$ws.Range("I136").Value = 2558592
$ws.Range("J136").Value = 2860.875
$ws.Range("K136").Value = 7675776
$ws.Range("L136").Value = 8582.625
$ws.Range("M136").Value = -7673226
$ws.Range("N136").Value = -13682.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 103.6
$ws.Range("J12").Value = 132.66667
$ws.Range("L12").Value = 398.00001
$ws.Range("N12").Value = -744.00001
$ws.Range("H129").Value = 44026.293
$ws.Range("I129").Value = 599.5
$ws.Range("J129").Value = 49816.535
$ws.Range("K129").Value = 1798.5
$ws.Range("L129").Value = 149449.605
$ws.Range("M129").Value = 3201.5
$ws.Range("N129").Value = -159449.605
$ws.Range("H131").Value = 793.88
$ws.Range("J131").Value = 811.086
$ws.Range("L131").Value = 2433.258
$ws.Range("N131").Value = -12513.258

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1808
$ws.Range("I97").Value = 1773.091
$ws.Range("K97").Value = 1773.091
$ws.Range("M97").Value = -1277.091
$ws.Range("H102").Value = 2119.375
$ws.Range("I102").Value = 2175.7273
$ws.Range("K102").Value = 2175.7273
$ws.Range("M102").Value = -553.7273
$ws.Range("H132").Value = 4276162
$ws.Range("I132").Value = 9617140
$ws.Range("K132").Value = 28851420
$ws.Range("M132").Value = -28848890
$ws.Range("H141").Value = 52248.75
$ws.Range("J141").Value = 52248.75
$ws.Range("L141").Value = 52248.75
$ws.Range("N141").Value = -62608.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3150
$ws.Range("J22").Value = 2350
$ws.Range("L22").Value = 2350
$ws.Range("N22").Value = -2940
$ws.Range("H27").Value = 3150
$ws.Range("J27").Value = 2350
$ws.Range("L27").Value = 2350
$ws.Range("N27").Value = -2564
$ws.Range("H46").Value = 2988
$ws.Range("I46").Value = 2200
$ws.Range("J46").Value = 3618.4
$ws.Range("K46").Value = 2200
$ws.Range("L46").Value = 3618.4
$ws.Range("M46").Value = -2012
$ws.Range("N46").Value = -3994.4
$ws.Range("H60").Value = 19999.5
$ws.Range("J60").Value = 19999.5
$ws.Range("L60").Value = 19999.5
$ws.Range("N60").Value = -21017.5
$ws.Range("H132").Value = 2896.5186
$ws.Range("I132").Value = 2211
$ws.Range("K132").Value = 6633
$ws.Range("M132").Value = -4103
$ws.Range("H136").Value = 3791.739
$ws.Range("I136").Value = 2814
$ws.Range("K136").Value = 8442
$ws.Range("M136").Value = -5892

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2860
$ws.Range("H132").Value = 1273.76
$ws.Range("I132").Value = 878.3171
$ws.Range("K132").Value = 2634.9513
$ws.Range("M132").Value = -104.9512999999997
$ws.Range("H141").Value = 76569.92999999999
$ws.Range("J141").Value = 76569.92999999999
$ws.Range("L141").Value = 76569.92999999999
$ws.Range("N141").Value = -86929.92999999999
